$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 73
$ws.Cells.Item(73, 1).Value = "Kartikey Gupta"
$ws.Cells.Item(73, 2).NumberFormat = "@"
$ws.Cells.Item(73, 2).Value = "2025-09-17"
$ws.Cells.Item(73, 2).ClearFormats()
$ws.Cells.Item(73, 3).Value = "20:32:45"

# Row 74
$ws.Cells.Item(74, 1).Value = "Kartikey Gupta"
$ws.Cells.Item(74, 2).NumberFormat = "@"
$ws.Cells.Item(74, 2).Value = "2025-09-17"
$ws.Cells.Item(74, 2).ClearFormats()
$ws.Cells.Item(74, 3).Value = "20:42:40"
